# Fix capitalization typo: "A model for" -> "A Model for"
# occurs on the title slide (slide 1) and on the agenda/recap slide (slide 2).

$p = $ppt.ActivePresentation

$old = "A model for Predicting Customer Retention in Telecom"
$new = "A Model for Predicting Customer Retention in Telecom"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if (-not $shape.HasTextFrame) { continue }
        $tf = $shape.TextFrame
        if (-not $tf.HasText) { continue }
        $tr = $tf.TextRange
        $full = $tr.Text
        $idx = $full.IndexOf($old)
        while ($idx -ge 0) {
            # Replace just the matched substring so any sibling runs
            # (e.g. the "TELCO CHURN PREDICTION MODEL" title run plus the
            # line break before this one) keep their own formatting.
            $sub = $tr.Characters($idx + 1, $old.Length)
            $sub.Text = $new
            $full = $tr.Text
            $idx = $full.IndexOf($old, $idx + $new.Length)
        }
    }
}
